$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = 1111143
$ws.Range("B29").Value = "SHIJIN MATHEW"
$ws.Range("C29").Value = 20
$ws.Range("D29").Value = "MALE"
$ws.Range("E29").Value = 4235252352
$ws.Range("F29").Value = "shijinmathew95@gmail.com"
$ws.Range("G29").Value = "RT-PCR"
$ws.Range("H29").Value = "NIHAL TIWARI"
$ws.Range("I29").Value = 500
$ws.Range("J29").Value = "23/03/2022"
$ws.Range("K29").Value = "04:56 PM"
